$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Three new log entries appended on 2025-08-11 (p2p and first mile runs).
# Columns: A=Timestamp, B=Quote ID, C=Status, D=Message
$timestamps = @("2025-08-11 12:32:57", "2025-08-11 12:41:52", "2025-08-11 12:57:18")
$quoteIds   = @("CMM0246LCL0044", "CMM0246LCL0039", "CMM0246LCL0044")
$statuses   = @("Success", "Success", "Success")
$messages   = @(
    "Origin: Nhava Sheva, India(INNSA); Shipment Scope: Port-to-Door; Entity: Arora Foods",
    "Origin: 110020, Okhla Industrial Estate, South East Delhi, Delhi, India; Shipment Scope: Door-to-Door; Entity: Arora Foods",
    "Origin: Nhava Sheva, India(INNSA); Shipment Scope: Port-to-Door; Entity: Arora Foods"
)

$startRow = 3

for ($i = 0; $i -lt $timestamps.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $timestamps[$i]
    $ws.Cells.Item($r, 2).Value = $quoteIds[$i]
    $ws.Cells.Item($r, 3).Value = $statuses[$i]
    $ws.Cells.Item($r, 4).Value = $messages[$i]
}

$wb.Save()
